# Apply the edit described by the diff:
#  1. Rename sheet "Interventions coverages" -> "Interventions cost and coverage"
#  2. Remove the four birth-order/time worksheets:
#       "RR birth by type", "birth distribution", "time between births", "RR birth by time"

$wb = $excel.ActiveWorkbook

# Suppress the "are you sure you want to delete" confirmation dialog.
$excel.DisplayAlerts = $false

# 1. Rename the "Interventions coverages" sheet.
$wb.Worksheets.Item("Interventions coverages").Name = "Interventions cost and coverage"

# 2. Delete the no-longer-needed worksheets.
$sheetsToRemove = @("RR birth by type", "birth distribution", "time between births", "RR birth by time")
foreach ($sheetName in $sheetsToRemove) {
    $wb.Worksheets.Item($sheetName).Delete()
}

$excel.DisplayAlerts = $true
